$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-18
# from 45208 (2023-10-09) to 45212 (2023-10-13)
$ws.Range("C2:C18").Value = 45212
